$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. The "Ticks" column headers (B1, E1, H1) are renamed to "Number of ticks"
#    (the same text already used by K1/N1). Once every reference to the old
#    "Ticks" shared string is gone, it drops out of the workbook on save.
$ws.Range("B1").Value = "Number of ticks"
$ws.Range("E1").Value = "Number of ticks"
$ws.Range("H1").Value = "Number of ticks"

# Those headers also lose their distinct "2-decimal number format" look and
# pick up the same plain bold/bordered style already used by every other
# header cell (copy the formatting from A1, which already has it).
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 2. The tick-count data columns (B, E, H) switch from a 2-decimal number
#    format to a plain integer format.
$ws.Range("B2:B7").NumberFormat = "0"
$ws.Range("E2:E7").NumberFormat = "0"
$ws.Range("H2:H7").NumberFormat = "0"

# 3. Those columns are widened a bit to fit the longer header text.
$ws.Columns("B").ColumnWidth = 14.67
$ws.Columns("E").ColumnWidth = 15.17
$ws.Columns("H").ColumnWidth = 15.17

# 4. The saved selection moves from H13 to K14.
$ws.Range("K14").Select()
